# "Registrar Oferta del fabricante"
# Adds three new comment rows to the errors/comments sheet (sheet2 - "Errores")
# describing issues found/fixed around registering an offer ("oferta") for a
# manufacturer's auction, and shifts the two subsequent entries
# (AuctionManagementBean / BussinessException) down to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert 3 blank rows before row 9, pushing existing rows 9-12 down to 12-15.
$ws.Rows("9:11").Insert()

# Populate the new comment rows (column C only, continuing under the same
# "TransactManager" grouping started in row 2).
$ws.Range("C9").Value = "El unico criterio actual es el mejor precio, el metodo de dar ganador no tiene sentido, el mejor se asigna cada vez que se registra una nueva oferta"
$ws.Range("C10").Value = "Al registrar la oferta no se estaba asignando a la subasta correspondiente"
$ws.Range("C11").Value = "El metodo de registrar oferta estaba en el bean pero no en el web service y por tanto nunca era empleado"

# Match the saved selection state from the authored workbook.
$ws.Range("C12").Select()
